$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 currently holds "AUG-DERMA SILICON GEL 15 GM" and row 8 holds "F B12 20PIECES".
# Overwrite them in place with the data for "ORGASOL LIGHT CREAM" and "PRISBRINA CAPS"
# (these rows keep their own row heights/styles).
$ws.Cells.Item(7, 3).Value2 = "ORGASOL LIGHT CREAM"
$ws.Cells.Item(7, 8).Value2 = "0:0"
$ws.Cells.Item(7, 12).Value2 = 0
$ws.Cells.Item(7, 14).Value2 = "130.00"

# P column (selling price) is numeric-formatted ("0.00"), but the source value is the
# text "130.0000" - force it to stay text without changing the cell's stored style.
$pCell7 = $ws.Cells.Item(7, 16)
$pCell7.NumberFormat = "@"
$pCell7.Value2 = "130.0000"
$pCell7.NumberFormat = "0.00"

$ws.Cells.Item(8, 3).Value2 = "PRISBRINA  CAPS"
$ws.Cells.Item(8, 8).Value2 = "-1:-1"
$ws.Cells.Item(8, 12).Value2 = 0
$ws.Cells.Item(8, 14).Value2 = "150.00"

$pCell8 = $ws.Cells.Item(8, 16)
$pCell8.NumberFormat = "@"
$pCell8.Value2 = "150.0000"
$pCell8.NumberFormat = "0.00"

# The old rows for "INJOCEPH 1000MG VIAL", the old "ORGASOL LIGHT CREAM" row, and the
# old "PRISBRINA CAPS" row are no longer needed now that their data moved up into rows 7-8.
$ws.Rows("9:11").Delete()

# Update the total (sum of selling prices) for the remaining rows: 130 + 150 = 280
$ws.Cells.Item(9, 16).Value2 = 280

# Update the generated timestamp in the footer
$ws.Cells.Item(10, 1).Value2 = "Saturday, 24 May, 2025 9:30 AM"
